$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" note text with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = [string]$wsHoja1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 4.1 = 15748.2 pesos"), "1000 Bs = 4.14 = 15953.66 pesos"
$newText = $newText -replace [regex]::Escape("15748.2 pesos = 4.07 = 935.7 Bs"), "15953.66 pesos = 4.13 = 950.61 Bs"
$wsHoja1.Range("A1").Value2 = $newText

# --- Update the "tasas" sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 241.7
$wsTasas.Range("O10").Value2 = 3856
$wsTasas.Range("N12").Value2 = 3860
